$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.384.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.55%  "
$ws.Range("D3").Value = "'1.776.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.86%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'307.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").Value = "'0.4232"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.60%  "
$ws.Range("D8").Value = "'0.3598"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.10%  "
$ws.Range("D9").Value = "'0.07132"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.67%  "
$ws.Range("D10").Value = "'0.8361"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("D11").Value = "'20.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.29%  "
$ws.Range("D12").Value = "'1.775.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.97%  "
$ws.Range("D13").Value = "'6.447"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.53%  "
$ws.Range("D14").Value = "'5.240"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").Value = "'0.06859"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("D16").Value = "'1.006"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("D17").Value = "'79.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").Value = "'0.000008631"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").Value = "'1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").Value = "'14.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").Value = "'26.388.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.42%  "
$ws.Range("D22").Value = "'5.075"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.95%  "
$ws.Range("D23").Value = "'10.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.07%  "
$ws.Range("D24").Value = "'1.999.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.88%  "
$ws.Range("D25").Value = "'152.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("D26").Value = "'1.814"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.47%  "
$ws.Range("D27").Value = "'17.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").Value = "'5.061"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("D29").Value = "'114.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.50%  "
$ws.Range("D30").Value = "'1.832"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.19%  "
$ws.Range("D31").Value = "'0.08821"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("D32").Value = "'0.7254"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("E33").Value = "  +5.01%  "
$ws.Range("D34").Value = "'4.321"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("E36").Value = "  -7.00%  "
$ws.Range("D37").Value = "'1.093"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.42%  "
$ws.Range("D38").Value = "'0.05127"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.00%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.1606"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.4906"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("D42").Value = "'2.600"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.69%  "
$ws.Range("D43").Value = "'6.331"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.73%  "
$ws.Range("D44").Value = "'7.962"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("D45").Value = "'104.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'10.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'1.000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("E48").Value = "  +3.79%  "
$ws.Range("D49").Value = "'0.06171"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.31%  "
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").Value = "'1.724"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.42%  "
